$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 3
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.63
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.83
$ws.Range("X7").Value = 11
$ws.Range("AI7").Value = 11
$ws.Range("AK7").Value = 26
$ws.Range("AT7").Value = 2.63
$ws.Range("AX7").Value = 17
$ws.Range("BB7").Value = 201
